# This script applies a weekly refresh to the "Espinaca" (Feria Lagunitas
# de Puerto Montt) sheet: a new week of data is inserted at row 12, which
# pushes all the previously-recorded weeks down by one row (old rows
# 12-48 become rows 13-49), and a brand-new date/price observation is
# written into the now-vacated row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 12..48 down into rows 13..49, working from the bottom up so
# that a row is never overwritten before it has been copied.  A full
# A:R row copy keeps every column (including the ones that are constant
# for this sheet) consistent with the source row.
for ($r = 48; $r -ge 12; $r--) {
    $target = $r + 1
    $src = $ws.Range("A" + $r + ":R" + $r)
    $dst = $ws.Range("A" + $target + ":R" + $target)
    $dst.Value = $src.Value2
}

# Row 49 did not exist before this edit, so its date cell needs to pick
# up the same number format used by the rest of column D.
$ws.Range("D49").NumberFormat = $ws.Range("D48").NumberFormat

# Write the new week's observation into row 12 (volume/prices are
# unchanged from the prior entry that used to sit in row 12; only the
# date is new).
$ws.Range("D12").Value = 44831
$ws.Range("J12").Value = 35
$ws.Range("K12").Value = 13000
$ws.Range("L12").Value = 13000
$ws.Range("M12").Value = 13000
$ws.Range("P12").Value = 1300
